$wb = $excel.ActiveWorkbook

# ALC row 2
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 3999.3635
$ws.Range("I2").Value = 6699.3335
$ws.Range("J2").Value = 759.4
$ws.Range("K2").Value = 6699.3335
$ws.Range("L2").Value = 759.4
$ws.Range("M2").Value = -6586.3335
$ws.Range("N2").Value = -985.4

# ALC row 19
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1739.5
$ws.Range("J19").Value = 1930.7142
$ws.Range("L19").Value = 1930.7142
$ws.Range("N19").Value = -2280.7142

# ALC row 38
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1057.6666
$ws.Range("J38").Value = 8999
$ws.Range("L38").Value = 26997
$ws.Range("N38").Value = -27741

# ALC row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1268.6957
$ws.Range("I98").Value = 1193.6842
$ws.Range("K98").Value = 1193.6842
$ws.Range("M98").Value = 304.3158000000001

# ALC row 106
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 25001428
$ws.Range("I106").Value = 27779180
$ws.Range("K106").Value = 27779180
$ws.Range("M106").Value = -27778549

# ALC row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1268.6957
$ws.Range("I122").Value = 1193.6842
$ws.Range("K122").Value = 3581.0526
$ws.Range("M122").Value = -1131.0526

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 7282082
$ws.Range("I45").Value = 11759745
$ws.Range("J45").Value = 5879.75
$ws.Range("K45").Value = 11759745
$ws.Range("L45").Value = 5879.75
$ws.Range("M45").Value = -11759368
$ws.Range("N45").Value = -6633.75

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3780.8948
$ws.Range("I61").Value = 3889.25
$ws.Range("J61").Value = 3477.5
$ws.Range("K61").Value = 3889.25
$ws.Range("L61").Value = 3477.5
$ws.Range("M61").Value = -3677.25
$ws.Range("N61").Value = -3901.5

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 50661.54
$ws.Range("I74").Value = 46086.75
$ws.Range("K74").Value = 46086.75
$ws.Range("M74").Value = -45212.75

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 50661.54
$ws.Range("I77").Value = 46086.75
$ws.Range("K77").Value = 230433.75
$ws.Range("M77").Value = -226065.75

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2164153.5
$ws.Range("I122").Value = 2508412.5
$ws.Range("K122").Value = 7525237.5
$ws.Range("M122").Value = -7522787.5

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3780.8948
$ws.Range("I136").Value = 3889.25
$ws.Range("J136").Value = 3477.5
$ws.Range("K136").Value = 11667.75
$ws.Range("L136").Value = 10432.5
$ws.Range("M136").Value = -9117.75
$ws.Range("N136").Value = -15532.5

# BSM row 80
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 572.5
$ws.Range("I80").Value = 596.5
$ws.Range("J80").Value = 568.5
$ws.Range("K80").Value = 596.5
$ws.Range("L80").Value = 568.5
$ws.Range("M80").Value = 401.5
$ws.Range("N80").Value = -2564.5

# BSM row 83
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 572.5
$ws.Range("I83").Value = 596.5
$ws.Range("J83").Value = 568.5
$ws.Range("K83").Value = 2982.5
$ws.Range("L83").Value = 2842.5
$ws.Range("M83").Value = 2009.5
$ws.Range("N83").Value = -12826.5

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 6506.095
$ws.Range("I134").Value = 1419.2941
$ws.Range("K134").Value = 4257.8823
$ws.Range("M134").Value = -1722.8823

# CRP row 26
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H26").Value = 32000
$ws.Range("J26").Value = 32000
$ws.Range("L26").Value = 32000
$ws.Range("N26").Value = -32574

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 62721.066
$ws.Range("I132").Value = 41410.28
$ws.Range("J132").Value = 151516
$ws.Range("K132").Value = 124230.84
$ws.Range("L132").Value = 454548
$ws.Range("M132").Value = -121700.84
$ws.Range("N132").Value = -459608

# CRP row 141
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 104752.05
$ws.Range("J141").Value = 109045.5
$ws.Range("L141").Value = 109045.5
$ws.Range("N141").Value = -119405.5

# CUL row 2
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 221.04762
$ws.Range("I2").Value = 67.052635
$ws.Range("J2").Value = 348.26086
$ws.Range("K2").Value = 402.3158099999999
$ws.Range("L2").Value = 2089.56516
$ws.Range("M2").Value = -289.3158099999999
$ws.Range("N2").Value = -2315.56516

# CUL row 12
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 150046.33
$ws.Range("I12").Value = 222494.5
$ws.Range("J12").Value = 5150
$ws.Range("K12").Value = 667483.5
$ws.Range("L12").Value = 15450
$ws.Range("M12").Value = -667310.5
$ws.Range("N12").Value = -15796

# CUL row 33
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 6749.7334
$ws.Range("I33").Value = 60.81818
$ws.Range("K33").Value = 364.90908
$ws.Range("M33").Value = -81.90908000000002

# CUL row 38
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 208.04762
$ws.Range("I38").Value = 197.8
$ws.Range("J38").Value = 233.66667
$ws.Range("K38").Value = 593.4000000000001
$ws.Range("L38").Value = 701.00001
$ws.Range("M38").Value = -246.4000000000001
$ws.Range("N38").Value = -1395.00001

# CUL row 44
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 2896.875
$ws.Range("I44").Value = 1235
$ws.Range("K44").Value = 3705
$ws.Range("M44").Value = -3307

# CUL row 63
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 2999.5
$ws.Range("J63").Value = 2000
$ws.Range("L63").Value = 6000
$ws.Range("N63").Value = -7498

# CUL row 64
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 1162.5
$ws.Range("I64").Value = 350
$ws.Range("J64").Value = 1433.3334
$ws.Range("K64").Value = 1050
$ws.Range("L64").Value = 4300.0002
$ws.Range("M64").Value = -780
$ws.Range("N64").Value = -4840.0002

# CUL row 66
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H66").Value = 2999.5
$ws.Range("J66").Value = 2000
$ws.Range("L66").Value = 18000
$ws.Range("N66").Value = -25488

# CUL row 67
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H67").Value = 1162.5
$ws.Range("I67").Value = 350
$ws.Range("J67").Value = 1433.3334
$ws.Range("K67").Value = 1050
$ws.Range("L67").Value = 4300.0002
$ws.Range("M67").Value = -114
$ws.Range("N67").Value = -6172.0002

# GSM row 55
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 7515
$ws.Range("I55").Value = 1030
$ws.Range("J55").Value = 14000
$ws.Range("K55").Value = 1030
$ws.Range("L55").Value = 14000
$ws.Range("M55").Value = -703
$ws.Range("N55").Value = -14654

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 19511016
$ws.Range("I80").Value = 30858202
$ws.Range("J80").Value = 220803.8
$ws.Range("K80").Value = 30858202
$ws.Range("L80").Value = 220803.8
$ws.Range("M80").Value = -30857204
$ws.Range("N80").Value = -222799.8

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 19511016
$ws.Range("I83").Value = 30858202
$ws.Range("J83").Value = 220803.8
$ws.Range("K83").Value = 154291010
$ws.Range("L83").Value = 1104019
$ws.Range("M83").Value = -154286018
$ws.Range("N83").Value = -1114003

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4956631.5
$ws.Range("I126").Value = 3499005.8
$ws.Range("K126").Value = 10497017.4
$ws.Range("M126").Value = -10494547.4

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6718.64
$ws.Range("I7").Value = 4763.3125
$ws.Range("J7").Value = 10194.777
$ws.Range("K7").Value = 4763.3125
$ws.Range("L7").Value = 10194.777
$ws.Range("M7").Value = -4651.3125
$ws.Range("N7").Value = -10418.777

# LTW row 82
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 42223816
$ws.Range("I82").Value = 95961350
$ws.Range("J82").Value = 1472.2858
$ws.Range("K82").Value = 95961350
$ws.Range("L82").Value = 1472.2858
$ws.Range("M82").Value = -95960989
$ws.Range("N82").Value = -2194.2858

# LTW row 85
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 42223816
$ws.Range("I85").Value = 95961350
$ws.Range("J85").Value = 1472.2858
$ws.Range("K85").Value = 95961350
$ws.Range("L85").Value = 1472.2858
$ws.Range("M85").Value = -95960102
$ws.Range("N85").Value = -3968.2858

# LTW row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 27796834
$ws.Range("I93").Value = 37039780
$ws.Range("J93").Value = 67994.336
$ws.Range("K93").Value = 37039780
$ws.Range("L93").Value = 67994.336
$ws.Range("M93").Value = -37038532
$ws.Range("N93").Value = -70490.336

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 6718.64
$ws.Range("I126").Value = 4763.3125
$ws.Range("J126").Value = 10194.777
$ws.Range("K126").Value = 14289.9375
$ws.Range("L126").Value = 30584.331
$ws.Range("M126").Value = -11819.9375
$ws.Range("N126").Value = -35524.331

# WVR row 4
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 35349.832
$ws.Range("J4").Value = 41419.8
$ws.Range("L4").Value = 41419.8
$ws.Range("N4").Value = -41645.8
